$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.26000000000001

$ws.Range("C7").Value = -12.6506
$ws.Range("D7").Value = -7.372100000000001

$ws.Range("D15").Value = -8.469199999999995

$ws.Range("C16").Value = -14.3757
$ws.Range("E16").Value = 16.2139

$ws.Range("E19").Value = 16.47099999999999

$ws.Range("D21").Value = -8.898899999999996

$ws.Range("D22").Value = -8.173800000000005

$ws.Range("D23").Value = -7.291900000000001

$ws.Range("C28").Value = -12.6984

$ws.Range("C29").Value = -11.12940000000001

$ws.Range("C32").Value = -13.2277

$ws.Range("D34").Value = -7.963000000000002

$ws.Range("E36").Value = 16.19380000000001

$ws.Range("C40").Value = -12.4659

$ws.Range("D43").Value = -8.346900000000003

$ws.Range("D45").Value = -7.852400000000001

$ws.Range("E46").Value = 17.13909999999998

$ws.Range("D50").Value = -8.499199999999995
$ws.Range("E50").Value = 16.63099999999999

$ws.Range("D51").Value = -7.590799999999998

$ws.Range("C52").Value = -11.0812

$ws.Range("C57").Value = -14.36729999999999

$ws.Range("C66").Value = -12.1515
$ws.Range("D66").Value = -7.794900000000001

$ws.Range("D67").Value = -6.669499999999998

$ws.Range("D79").Value = -6.195100000000005

$ws.Range("D84").Value = -8.698400000000003

$ws.Range("D92").Value = -6.433200000000004

$ws.Range("E95").Value = 17.97040000000002

$ws.Range("D97").Value = -8.259799999999997
$ws.Range("E97").Value = 16.34569999999999

$ws.Range("C100").Value = -12.6215
